$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data as per latest scrape (GitHub Actions run)
# Price column (D) values are forced to text format to preserve exact
# textual representation (e.g. "1.000", "28.226.64") instead of being
# auto-converted to numeric values by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.226.64"
$ws.Range("E2").Value = "  +2.42%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.870.23"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.60%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "337.09"
$ws.Range("E5").Value = "  +1.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4699"
$ws.Range("E7").Value = "  +1.65%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3926"
$ws.Range("E8").Value = "  +1.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.10"
$ws.Range("E9").Value = "  +2.58%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07991"
$ws.Range("E10").Value = "  +1.22%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.011"
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.71"
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.871.19"
$ws.Range("E13").Value = "  +0.47%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.994"
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.264"
$ws.Range("E15").Value = "  +2.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.31"
$ws.Range("E16").Value = "  +3.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001041"
$ws.Range("E18").Value = "  +0.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06610"
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.73"
$ws.Range("E20").Value = "  +3.84%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9998"
$ws.Range("E21").Value = "  -0.66%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "28.223.97"
$ws.Range("E22").Value = "  +2.48%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.441"
$ws.Range("E23").Value = "  +1.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.08"
$ws.Range("E24").Value = "  +1.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.295"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.076.04"
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.02"
$ws.Range("E27").Value = "  +0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.04"
$ws.Range("E28").Value = "  +2.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.145"
$ws.Range("E29").Value = "  +1.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.487"
$ws.Range("E30").Value = "  +1.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "119.88"
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9767"
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09468"
$ws.Range("E33").Value = "  +0.99%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.576"
$ws.Range("E34").Value = "  -0.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.377"
$ws.Range("E35").Value = "  +2.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.349"
$ws.Range("E36").Value = "  +1.47%  "
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06106"
$ws.Range("E37").Value = "  +1.34%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02260"
$ws.Range("E38").Value = "  +1.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.439"
$ws.Range("E39").Value = "  +2.03%  "
$ws.Range("E40").Value = "  -0.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5969"
$ws.Range("E41").Value = "  +1.54%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.0000"
$ws.Range("E42").Value = "  -0.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1878"
$ws.Range("E43").Value = "  +0.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.39"
$ws.Range("E44").Value = "  +1.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.293"
$ws.Range("E45").Value = "  +4.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5625"
$ws.Range("E46").Value = "  +0.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.13"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.974"
$ws.Range("E48").Value = "  +3.88%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06868"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "111.28"
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.972"
$ws.Range("E51").Value = "  +11.70%  "